# Auto-generated Excel COM-interop script to apply the diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.63
# Row 4
$ws.Range("I4").Value = 2.01
# Row 5
$ws.Range("AD5").Value = 8
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 23
$ws.Range("AN5").Value = 3.2
$ws.Range("AP5").Value = 29
$ws.Range("G5").Value = 1.51
$ws.Range("I5").Value = 7
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 2.38
$ws.Range("U5").Value = 2.75
$ws.Range("V5").Value = 1.4
$ws.Range("X5").Value = 5.5
$ws.Range("Z5").Value = 10
# Row 6
$ws.Range("AA6").Value = 51
$ws.Range("AC6").Value = 6.5
$ws.Range("AJ6").Value = 9.5
$ws.Range("AM6").Value = 41
$ws.Range("AQ6").Value = 126
$ws.Range("AV6").Value = 81
$ws.Range("AW6").Value = 3.5
$ws.Range("G6").Value = 5.25
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 1.71
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("W6").Value = 10
$ws.Range("Y6").Value = 19
# Row 7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
# Row 10
$ws.Range("AC10").Value = 7
$ws.Range("AJ10").Value = 19
$ws.Range("AR10").Value = 67
$ws.Range("AT10").Value = 2.5
$ws.Range("AU10").Value = 9.5
$ws.Range("BA10").Value = 151
$ws.Range("H10").Value = 3.2
$ws.Range("J10").Value = 2.5
$ws.Range("K10").Value = 2.05
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.63
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("U10").Value = 2.2
$ws.Range("V10").Value = 1.62
$ws.Range("W10").Value = 5.5
$ws.Range("X10").Value = 7
# Row 12
$ws.Range("AL12").Value = 41
$ws.Range("AZ12").Value = 101
$ws.Range("I12").Value = 5
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 6.5
$ws.Range("U12").Value = 2.2
$ws.Range("V12").Value = 1.62
$ws.Range("Y12").Value = 9
# Row 13
$ws.Range("Y13").Value = 12
